$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) store numeric/percent-looking
# values as literal text in the source data. Force each target cell
# to Text format first so Excel does not auto-convert the assigned
# string into a Number/Percentage value.
foreach ($addr in @("D2", "D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "E2", "E3", "E4", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E17", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E28", "E40", "E41", "E42", "E43", "E44", "E45", "E46", "E47", "E48", "E49", "E50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '245.65'
$ws.Range("E2").Value = '1.26%'
$ws.Range("E3").Value = '1.97%'
$ws.Range("D4").Value = '5.161'
$ws.Range("E4").Value = '1.00%'
$ws.Range("D5").Value = '0.05750'
$ws.Range("D6").Value = '6.566'
$ws.Range("E6").Value = '1.14%'
$ws.Range("D7").Value = '0.8573'
$ws.Range("E7").Value = '3.84%'
$ws.Range("D8").Value = '0.8553'
$ws.Range("E8").Value = '-1.53%'
$ws.Range("E9").Value = '2.20%'
$ws.Range("D10").Value = '0.07021'
$ws.Range("E10").Value = '1.97%'
$ws.Range("D11").Value = '0.03045'
$ws.Range("E11").Value = '6.51%'
$ws.Range("D12").Value = '0.09358'
$ws.Range("E12").Value = '-0.12%'
$ws.Range("D13").Value = '0.001535'
$ws.Range("E13").Value = '1.43%'
$ws.Range("D14").Value = '0.0005979'
$ws.Range("E14").Value = '-0.51%'
$ws.Range("D15").Value = '0.005988'
$ws.Range("E15").Value = '-1.59%'
$ws.Range("D16").Value = '3.495'
$ws.Range("E16").Value = '-0.78%'
$ws.Range("D17").Value = '3.117'
$ws.Range("E17").Value = '3.16%'
$ws.Range("D18").Value = '2.218'
$ws.Range("E18").Value = '0.15%'
$ws.Range("D19").Value = '0.3200'
$ws.Range("E19").Value = '1.61%'
$ws.Range("D20").Value = '0.03282'
$ws.Range("E20").Value = '1.15%'
$ws.Range("D21").Value = '0.1275'
$ws.Range("E21").Value = '0.18%'
$ws.Range("D22").Value = '3.510'
$ws.Range("E22").Value = '-2.92%'
$ws.Range("D23").Value = '0.04159'
$ws.Range("E23").Value = '-0.34%'
$ws.Range("E24").Value = '0.42%'
$ws.Range("D25").Value = '0.001226'
$ws.Range("E25").Value = '1.28%'
$ws.Range("D26").Value = '0.004134'
$ws.Range("E26").Value = '-6.94%'
$ws.Range("E27").Value = '2.55%'
$ws.Range("E28").Value = '3.18%'
$ws.Range("D40").Value = '0.03724'
$ws.Range("E40").Value = '0.56%'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1071'
$ws.Range("E41").Value = '1.62%'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '0.002460'
$ws.Range("E42").Value = '6.42%'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '0.003500'
$ws.Range("E43").Value = '-39.54%'
$ws.Range("D44").Value = '0.009323'
$ws.Range("E44").Value = '-2.47%'
$ws.Range("D45").Value = '0.00005272'
$ws.Range("E45").Value = '3.63%'
$ws.Range("E46").Value = '-0.02%'
$ws.Range("D47").Value = '0.05799'
$ws.Range("E47").Value = '-51.67%'
$ws.Range("D48").Value = '0.002448'
$ws.Range("E48").Value = '-1.39%'
$ws.Range("E49").Value = '-0.02%'
$ws.Range("E50").Value = '-0.02%'
